$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 00:57"

# Row 4
$ws.Range("B4").Value = 7598590
$ws.Range("C4").Value = 46669
$ws.Range("D4").Value = 4813737
$ws.Range("E4").Value = 2570594
$ws.Range("G4").Value = 737
$ws.Range("H4").Value = 214259

# Row 6
$ws.Range("B6").Value = 4906833
$ws.Range("C6").Value = 24602
$ws.Range("D6").Value = 4248574
$ws.Range("E6").Value = 512272
$ws.Range("G6").Value = 556
$ws.Range("H6").Value = 145987

# Row 8
$ws.Range("B8").Value = 848147
$ws.Range("C8").Value = 6616
$ws.Range("D8").Value = 757801
$ws.Range("E8").Value = 63790
$ws.Range("G8").Value = 159
$ws.Range("H8").Value = 26556

# Row 9
$ws.Range("B9").Value = 824985
$ws.Range("C9").Value = 3421
$ws.Range("D9").Value = 700868
$ws.Range("E9").Value = 91452
$ws.Range("G9").Value = 56
$ws.Range("H9").Value = 32665

# Row 11
$ws.Range("B11").Value = 790818
$ws.Range("C11").Value = 11129
$ws.Range("D11").Value = 626114
$ws.Range("E11").Value = 143909
$ws.Range("G11").Value = 196
$ws.Range("H11").Value = 20795

# Row 37
$ws.Range("B37").Value = 114653
$ws.Range("C37").Value = 691
$ws.Range("D37").Value = 91195
$ws.Range("E37").Value = 21044
$ws.Range("G37").Value = 8
$ws.Range("H37").Value = 2414

# Row 41
$ws.Range("B41").Value = 103575
$ws.Range("C41").Value = 109
$ws.Range("D41").Value = 97274
$ws.Range("E41").Value = 331
$ws.Range("G41").Value = 14
$ws.Range("H41").Value = 5970

# Row 57
$ws.Range("B57").Value = 72310
$ws.Range("C57").Value = 507
$ws.Range("D57").Value = 66813
$ws.Range("E57").Value = 5239
$ws.Range("G57").Value = 3
$ws.Range("H57").Value = 258

# Row 58
$ws.Range("B58").Value = 59287
$ws.Range("C58").Value = 160
$ws.Range("D58").Value = 50718
$ws.Range("E58").Value = 7456
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 1113

# Row 83
$ws.Range("B83").Value = 21518
$ws.Range("C83").Value = 182
$ws.Range("D83").Value = 14984
$ws.Range("E83").Value = 5693
$ws.Range("G83").Value = 3
$ws.Range("H83").Value = 841

# Row 93
$ws.Range("B93").Value = 14974
$ws.Range("C93").Value = 144
$ws.Range("D93").Value = 14094
$ws.Range("E93").Value = 547

# Row 109
$ws.Range("A109").Value = "Haiti"
$ws.Range("B109").Value = 8811
$ws.Range("C109").Value = 19
$ws.Range("D109").Value = 6949
$ws.Range("E109").Value = 1633
$ws.Range("H109").Value = 229

# Row 110
$ws.Range("A110").Value = "Luxemburgo"
$ws.Range("C110").Value = 88
$ws.Range("D110").Value = 7428
$ws.Range("E110").Value = 1244
$ws.Range("H110").Value = 125

# Row 111
$ws.Range("A111").Value = "Gabon"
$ws.Range("B111").Value = 8797
$ws.Range("D111").Value = 8067
$ws.Range("E111").Value = 676
$ws.Range("H111").Value = 54

# Row 121
$ws.Range("A121").Value = "Suazilandia"
$ws.Range("B121").Value = 5530
$ws.Range("C121").Value = 9
$ws.Range("D121").Value = 5076
$ws.Range("E121").Value = 343
$ws.Range("H121").Value = 111

# Row 122
$ws.Range("A122").Value = "Guadalupe"
$ws.Range("B122").Value = 5528
$ws.Range("D122").Value = 2199
$ws.Range("E122").Value = 3272
$ws.Range("H122").Value = 57

# Row 123
$ws.Range("B123").Value = 5418
$ws.Range("C123").Value = 1
$ws.Range("D123").Value = 5347

# Row 132
$ws.Range("B132").Value = 4845
$ws.Range("C132").Value = 16
$ws.Range("E132").Value = 2869

# Row 152
$ws.Range("B152").Value = 2259
$ws.Range("C152").Value = 7
$ws.Range("D152").Value = 1704

# Row 154
$ws.Range("A154").Value = "Uruguay"
$ws.Range("B154").Value = 2122
$ws.Range("C154").Value = 25
$ws.Range("D154").Value = 1831
$ws.Range("E154").Value = 243
$ws.Range("H154").Value = 48

# Row 155
$ws.Range("A155").Value = "Principado de Andorra"
$ws.Range("B155").Value = 2110
$ws.Range("D155").Value = 1540
$ws.Range("E155").Value = 517
$ws.Range("H155").Value = 53

# Row 207
$ws.Range("A207").Value = "Nueva Caledonia"

# Row 208
$ws.Range("A208").Value = "Santa Lucia"
